$p = $ppt.ActivePresentation

$s = $p.Slides.Item(6)
$shp = $s.Shapes.Item(2)
$tbl = $shp.Table
$tbl.ApplyStyle("{749C632C-FA79-43EC-8786-6A699CD546E7}")
